$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'40.902.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.88%  '

$ws.Range("D3").Value = "'2.219.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.73%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = "'229.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.42%  '

$ws.Range("E6").Value = '  +2.11%  '

$ws.Range("D7").Value = "'64.73"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.52%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  +2.16%  '

$ws.Range("D10").Value = "'0.0869"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.23%  '

$ws.Range("E11").Value = '  +0.41%  '

$ws.Range("D12").Value = "'2.547.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.69%  '

$ws.Range("D13").Value = "'15.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.17%  '

$ws.Range("D14").Value = "'22.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.05%  '

$ws.Range("D15").Value = "'0.825"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.03%  '

$ws.Range("D16").Value = "'5.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.35%  '

$ws.Range("D17").Value = "'2.216.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.94%  '

$ws.Range("D18").Value = "'40.779.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.82%  '

$ws.Range("D19").Value = "'74.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.13%  '

$ws.Range("D20").Value = "'0.0₃0905"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.00%  '

$ws.Range("D21").Value = "'6.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.20%  '

$ws.Range("D22").Value = "'253.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.54%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("D24").Value = "'2.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.95%  '

$ws.Range("E25").Value = '  -8.07%  '

$ws.Range("D26").Value = "'9.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.65%  '

$ws.Range("D27").Value = "'173.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.64%  '

$ws.Range("D28").Value = "'0.143"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.05%  '

$ws.Range("D29").Value = "'20.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.37%  '

$ws.Range("E30").Value = '  +2.64%  '

$ws.Range("D31").Value = "'2.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.86%  '

$ws.Range("E32").Value = '  +1.35%  '

$ws.Range("D33").Value = "'4.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.27%  '

$ws.Range("D34").Value = "'7.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.25%  '

$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("D36").Value = "'0.0633"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.27%  '

$ws.Range("D37").Value = "'3.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.16%  '

$ws.Range("D38").Value = "'2.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.00%  '

$ws.Range("E39").Value = '  -0.04%  '

$ws.Range("D40").Value = "'4.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +13.85%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = "'0.0233"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.44%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = "'8.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.99%  '

$ws.Range("D43").Value = "'101.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.43%  '

$ws.Range("D44").Value = "'1.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.49%  '

$ws.Range("D45").Value = "'1.513.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.80%  '

$ws.Range("D46").Value = "'17.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.14%  '

$ws.Range("E47").Value = '  +1.44%  '

$ws.Range("E48").Value = '  +2.10%  '

$ws.Range("E49").Value = '  +0.00%  '

$ws.Range("E50").Value = '  +38.76%  '

$ws.Range("D51").Value = "'9.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +11.39%  '
